# Update the Salesperson column (B2:B23) on Sheet1 from "First Last" to
# "Last, First" format. Assigning new literal text to each cell adds a
# fresh shared-string entry (Excel only reuses an existing shared string
# when the text matches exactly), which grows sharedStrings.xml's
# uniqueCount from 102 to 124 and repoints B2:B23 at the new entries -
# exactly mirroring the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "Byers, Jordan"
$ws.Range("B3").Value  = "Collins, Casey"
$ws.Range("B4").Value  = "Gellers, Alex"
$ws.Range("B5").Value  = "Jones, Morgan"
$ws.Range("B6").Value  = "Lowe, Taylor"
$ws.Range("B7").Value  = "Quinn, Sam"
$ws.Range("B8").Value  = "Smith, Riley"
$ws.Range("B9").Value  = "Thomas, Blake"
$ws.Range("B10").Value = "Reynolds, Dylan"
$ws.Range("B11").Value = "Smith, James"
$ws.Range("B12").Value = "Owen, Michelle"
$ws.Range("B13").Value = "Ortega, Raul"
$ws.Range("B14").Value = "Chen, Maria"
$ws.Range("B15").Value = "Olsen, Leslie"
$ws.Range("B16").Value = "Nelson, Jamika"
$ws.Range("B17").Value = "Williams, Desmond"
$ws.Range("B18").Value = "Cooley, Renee"
$ws.Range("B19").Value = "Johnson, Tim"
$ws.Range("B20").Value = "Huen, James"
$ws.Range("B21").Value = "Ryeo, Annie"
$ws.Range("B22").Value = "Morgan, Emma"
$ws.Range("B23").Value = "Jones, Jackson"

# Widen column B so the longer "Last, First" names aren't truncated.
$ws.Columns("B").ColumnWidth = 17.5

# The source sheet's selection pointed at B1; move it back to the sheet's
# top-left (A1) now that the edit is done, matching the tidied-up view
# state in the commit.
$ws.Range("A1").Select()
